$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2783.3333
$ws.Range("I18").Value = 2843.75
$ws.Range("J18").Value = 2300
$ws.Range("K18").Value = 2843.75
$ws.Range("L18").Value = 2300
$ws.Range("M18").Value = -2559.75
$ws.Range("N18").Value = -2868
$ws.Range("H21").Value = 6714.5
$ws.Range("I21").Value = 6143.2
$ws.Range("K21").Value = 6143.2
$ws.Range("M21").Value = -5675.2
$ws.Range("H23").Value = 6714.5
$ws.Range("I23").Value = 6143.2
$ws.Range("K23").Value = 6143.2
$ws.Range("M23").Value = -5909.2
$ws.Range("H29").Value = 1191
$ws.Range("I29").Value = 949.5
$ws.Range("J29").Value = 1311.75
$ws.Range("K29").Value = 2848.5
$ws.Range("L29").Value = 3935.25
$ws.Range("M29").Value = -2567.5
$ws.Range("N29").Value = -4497.25
$ws.Range("H38").Value = 296.125
$ws.Range("I38").Value = 296.125
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 888.375
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -516.375
$ws.Range("N38").ClearContents()
$ws.Range("H43").Value = 4666.4
$ws.Range("I43").Value = 3567
$ws.Range("K43").Value = 3567
$ws.Range("M43").Value = -3498
$ws.Range("H55").Value = 177.85715
$ws.Range("I55").Value = 189.2
$ws.Range("K55").Value = 189.2
$ws.Range("M55").Value = 24.80000000000001
$ws.Range("H58").Value = 3831.7778
$ws.Range("I58").Value = 219.8
$ws.Range("J58").Value = 8346.75
$ws.Range("K58").Value = 659.4000000000001
$ws.Range("L58").Value = 25040.25
$ws.Range("M58").Value = -509.4000000000001
$ws.Range("N58").Value = -25340.25
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H80").Value = 3018.647
$ws.Range("J80").Value = 3520.125
$ws.Range("L80").Value = 10560.375
$ws.Range("N80").Value = -12556.375
$ws.Range("H83").Value = 3018.647
$ws.Range("J83").Value = 3520.125
$ws.Range("L83").Value = 31681.125
$ws.Range("N83").Value = -41665.125
$ws.Range("H87").Value = 189999
$ws.Range("J87").Value = 189999
$ws.Range("L87").Value = 189999
$ws.Range("N87").Value = -192495
$ws.Range("H90").Value = 189999
$ws.Range("J90").Value = 189999
$ws.Range("L90").Value = 569997
$ws.Range("N90").Value = -582477
$ws.Range("H106").Value = 12056.625
$ws.Range("I106").Value = 13614.714
$ws.Range("K106").Value = 13614.714
$ws.Range("M106").Value = -12983.714
$ws.Range("H137").Value = 41667988
$ws.Range("I137").Value = 47620270
$ws.Range("K137").Value = 142860810
$ws.Range("M137").Value = -142858260
$ws.Range("H138").Value = 2431.2666
$ws.Range("I138").Value = 4547.625
$ws.Range("J138").Value = 1973.6757
$ws.Range("K138").Value = 13642.875
$ws.Range("L138").Value = 5921.0271
$ws.Range("M138").Value = -8502.875
$ws.Range("N138").Value = -16201.0271

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3790.389
$ws.Range("I32").Value = 3131.147
$ws.Range("J32").Value = 14997.5
$ws.Range("K32").Value = 3131.147
$ws.Range("L32").Value = 14997.5
$ws.Range("M32").Value = -2844.147
$ws.Range("N32").Value = -15571.5
$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H76").Value = 52236.8
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 52236.8
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13904.909
$ws.Range("I86").Value = 14572.723
$ws.Range("J86").Value = 10899.75
$ws.Range("K86").Value = 14572.723
$ws.Range("L86").Value = 10899.75
$ws.Range("M86").Value = -13449.723
$ws.Range("N86").Value = -13145.75
$ws.Range("H89").Value = 13904.909
$ws.Range("I89").Value = 14572.723
$ws.Range("J89").Value = 10899.75
$ws.Range("K89").Value = 72863.61500000001
$ws.Range("L89").Value = 54498.75
$ws.Range("M89").Value = -67247.61500000001
$ws.Range("N89").Value = -65730.75
$ws.Range("H94").Value = 2044.2858
$ws.Range("I94").Value = 1910
$ws.Range("J94").Value = 2223.3333
$ws.Range("K94").Value = 1910
$ws.Range("L94").Value = 2223.3333
$ws.Range("M94").Value = -1459
$ws.Range("N94").Value = -3125.3333

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 90114.8
$ws.Range("J20").Value = 90114.8
$ws.Range("L20").Value = 90114.8
$ws.Range("N20").Value = -90586.8
$ws.Range("H22").Value = 2777.3333
$ws.Range("I22").Value = 2230
$ws.Range("J22").Value = 3872
$ws.Range("K22").Value = 2230
$ws.Range("L22").Value = 3872
$ws.Range("M22").Value = -1880
$ws.Range("N22").Value = -4572
$ws.Range("H30").Value = 90114.8
$ws.Range("J30").Value = 90114.8
$ws.Range("L30").Value = 90114.8
$ws.Range("N30").Value = -90296.8
$ws.Range("H86").Value = 11141.5
$ws.Range("I86").Value = 12118.857
$ws.Range("K86").Value = 12118.857
$ws.Range("M86").Value = -10995.857
$ws.Range("H89").Value = 11141.5
$ws.Range("I89").Value = 12118.857
$ws.Range("K89").Value = 60594.285
$ws.Range("M89").Value = -54978.285
$ws.Range("H128").Value = 90114.8
$ws.Range("J128").Value = 90114.8
$ws.Range("L128").Value = 90114.8
$ws.Range("N128").Value = -100074.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2549
$ws.Range("I2").Value = 98
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 588
$ws.Range("L2").Value = 30000
$ws.Range("M2").Value = -475
$ws.Range("N2").Value = -30226
$ws.Range("H4").Value = 6700183.5
$ws.Range("I4").Value = 3787706.2
$ws.Range("K4").Value = 11363118.6
$ws.Range("M4").Value = -11363006.6
$ws.Range("H7").Value = 340067.34
$ws.Range("I7").Value = 340067.34
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1020202.02
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1020090.02
$ws.Range("N7").ClearContents()
$ws.Range("H12").Value = 73.125
$ws.Range("I12").Value = 34
$ws.Range("J12").Value = 78.71429000000001
$ws.Range("K12").Value = 102
$ws.Range("L12").Value = 236.14287
$ws.Range("M12").Value = 71
$ws.Range("N12").Value = -582.14287
$ws.Range("H26").Value = 985.8570999999999
$ws.Range("J26").Value = 1180.2
$ws.Range("L26").Value = 3540.6
$ws.Range("N26").Value = -4116.6
$ws.Range("H70").Value = 8000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 8000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 100000
$ws.Range("I80").Value = 100000
$ws.Range("K80").Value = 300000
$ws.Range("M80").Value = -299064
$ws.Range("H81").Value = 1890
$ws.Range("I81").Value = 1280
$ws.Range("K81").Value = 3840
$ws.Range("M81").Value = -2717
$ws.Range("H83").Value = 100000
$ws.Range("I83").Value = 100000
$ws.Range("K83").Value = 900000
$ws.Range("M83").Value = -895320
$ws.Range("H84").Value = 1890
$ws.Range("I84").Value = 1280
$ws.Range("K84").Value = 11520
$ws.Range("M84").Value = -5904
$ws.Range("H131").Value = 2171.6428
$ws.Range("J131").Value = 2031
$ws.Range("L131").Value = 6093
$ws.Range("N131").Value = -16173
$ws.Range("H138").Value = 8575.15
$ws.Range("I138").Value = 7041
$ws.Range("K138").Value = 21123
$ws.Range("M138").Value = -15983

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 45000
$ws.Range("J33").Value = 45000
$ws.Range("L33").Value = 45000
$ws.Range("N33").Value = -45504
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H80").Value = 5411.1665
$ws.Range("I80").Value = 2494.6667
$ws.Range("J80").Value = 8327.666999999999
$ws.Range("K80").Value = 2494.6667
$ws.Range("L80").Value = 8327.666999999999
$ws.Range("M80").Value = -1496.6667
$ws.Range("N80").Value = -10323.667
$ws.Range("H83").Value = 5411.1665
$ws.Range("I83").Value = 2494.6667
$ws.Range("J83").Value = 8327.666999999999
$ws.Range("K83").Value = 12473.3335
$ws.Range("L83").Value = 41638.335
$ws.Range("M83").Value = -7481.333500000001
$ws.Range("N83").Value = -51622.335
$ws.Range("H97").Value = 685.63635
$ws.Range("I97").Value = 464.2
$ws.Range("K97").Value = 464.2
$ws.Range("M97").Value = 31.80000000000001
$ws.Range("H128").Value = 106980
$ws.Range("J128").Value = 106980
$ws.Range("L128").Value = 106980
$ws.Range("N128").Value = -116940
$ws.Range("H132").Value = 13890958
$ws.Range("I132").Value = 1637.0625
$ws.Range("K132").Value = 4911.1875
$ws.Range("M132").Value = -2381.1875
$ws.Range("H133").Value = 78999.5
$ws.Range("J133").Value = 78999.5
$ws.Range("L133").Value = 78999.5
$ws.Range("N133").Value = -89119.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4312.6665
$ws.Range("I7").Value = 3628.6667
$ws.Range("K7").Value = 3628.6667
$ws.Range("M7").Value = -3516.6667
$ws.Range("H16").Value = 1916
$ws.Range("I16").Value = 1375
$ws.Range("J16").Value = 2998
$ws.Range("K16").Value = 1375
$ws.Range("L16").Value = 2998
$ws.Range("M16").Value = -1205
$ws.Range("N16").Value = -3338
$ws.Range("H42").Value = 39990
$ws.Range("J42").Value = 39990
$ws.Range("L42").Value = 39990
$ws.Range("N42").Value = -41116
$ws.Range("H43").Value = 29999
$ws.Range("J43").Value = 29999
$ws.Range("L43").Value = 29999
$ws.Range("N43").Value = -30385
$ws.Range("H49").Value = 39990
$ws.Range("J49").Value = 39990
$ws.Range("L49").Value = 39990
$ws.Range("N49").Value = -40284
$ws.Range("H55").Value = 839.0833
$ws.Range("I55").Value = 388.7143
$ws.Range("J55").Value = 1469.6
$ws.Range("K55").Value = 388.7143
$ws.Range("L55").Value = 1469.6
$ws.Range("M55").Value = -215.7143
$ws.Range("N55").Value = -1815.6
$ws.Range("H68").Value = 1024
$ws.Range("I68").Value = 1033.6666
$ws.Range("K68").Value = 1033.6666
$ws.Range("M68").Value = -284.6666
$ws.Range("H71").Value = 1024
$ws.Range("I71").Value = 1033.6666
$ws.Range("K71").Value = 5168.333000000001
$ws.Range("M71").Value = -1424.333000000001
$ws.Range("H76").Value = 14962.25
$ws.Range("J76").Value = 14962.25
$ws.Range("L76").Value = 14962.25
$ws.Range("N76").Value = -15638.25
$ws.Range("H79").Value = 14962.25
$ws.Range("J79").Value = 14962.25
$ws.Range("L79").Value = 14962.25
$ws.Range("N79").Value = -17302.25
$ws.Range("H82").Value = 1138
$ws.Range("I82").Value = 427
$ws.Range("J82").Value = 1331.909
$ws.Range("K82").Value = 427
$ws.Range("L82").Value = 1331.909
$ws.Range("M82").Value = -66
$ws.Range("N82").Value = -2053.909
$ws.Range("H85").Value = 1138
$ws.Range("I85").Value = 427
$ws.Range("J85").Value = 1331.909
$ws.Range("K85").Value = 427
$ws.Range("L85").Value = 1331.909
$ws.Range("M85").Value = 821
$ws.Range("N85").Value = -3827.909
$ws.Range("H93").Value = 851.7143
$ws.Range("I93").Value = 822.6
$ws.Range("K93").Value = 822.6
$ws.Range("M93").Value = 425.4
$ws.Range("H122").Value = 1541.6364
$ws.Range("I122").Value = 1541.6364
$ws.Range("K122").Value = 4624.9092
$ws.Range("M122").Value = -2174.9092
$ws.Range("H126").Value = 4312.6665
$ws.Range("I126").Value = 3628.6667
$ws.Range("K126").Value = 10886.0001
$ws.Range("M126").Value = -8416.000100000001
$ws.Range("H128").Value = 70255.5
$ws.Range("J128").Value = 70255.5
$ws.Range("L128").Value = 70255.5
$ws.Range("N128").Value = -80215.5
$ws.Range("H132").Value = 3550
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 47306
$ws.Range("J54").Value = 47306
$ws.Range("L54").Value = 47306
$ws.Range("N54").Value = -48346
$ws.Range("H107").Value = 1009.5
$ws.Range("I107").Value = 1081.7142
$ws.Range("J107").Value = 953.3333
$ws.Range("K107").Value = 3245.1426
$ws.Range("L107").Value = 2859.9999
$ws.Range("M107").Value = -1325.1426
$ws.Range("N107").Value = -6699.9999
$ws.Range("H130").Value = 24106.334
$ws.Range("J130").Value = 24106.334
$ws.Range("L130").Value = 24106.334
$ws.Range("N130").Value = -34146.334
$ws.Range("H132").Value = 142859940
$ws.Range("I132").Value = 3617.2
$ws.Range("K132").Value = 10851.6
$ws.Range("M132").Value = -8321.599999999999
